$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: correct the previously-blank formula row with real typed data ---
$ws.Range("A3").Value = "Test2"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 44196
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = $true

# --- Row 4: brand-new row appended below ---
$ws.Range("A4").Value = "Test3"
$ws.Range("D4").Value = 44196
$ws.Range("F4").Value = $true

# Reuse the existing DateTime/Currency number formats (style indexes) from
# row 2 instead of letting Excel mint new ones from a raw NumberFormat string.
$ws.Range("D2").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E2").Copy()
$ws.Range("E3:E4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Grow Table1 so its range (and autofilter) covers the new row too.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F4"))

$ws.Range("F4").Select()
